$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cell X1 "Unnamed: 23" with the same formatting as the other header cells
$ws.Range("X1").Value = "Unnamed: 23"
$ws.Range("W1").Copy()
$ws.Range("X1").PasteSpecial(-4122)
$ws.Range("X1").Value = "Unnamed: 23"
$excel.CutCopyMode = 0

# Add new data row 3 (keep Date column as plain text, not an auto-converted date)
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "09/12/2025"
$ws.Range("A3").ClearFormats()

$ws.Range("B3").Value = "Qwen2.5-32B-Instruct"
$ws.Range("C3").Value = 0.4953703703703703
$ws.Range("D3").Value = 0.3602693602693603
$ws.Range("E3").Value = 0.4171539961013645
$ws.Range("F3").Value = 0.5385039129371485
$ws.Range("G3").Value = 0.366767344502872
$ws.Range("H3").Value = 0.4254899499088304
$ws.Range("I3").Value = 0.546178414124122
$ws.Range("J3").Value = 0.3602693602693603
$ws.Range("K3").Value = 0.4254327582600828
$ws.Range("L3").Value = 0.5833333333333334
$ws.Range("M3").Value = 0.4256756756756757
$ws.Range("N3").Value = 0.4921875000000001
$ws.Range("O3").Value = 126
$ws.Range("P3").Value = 90
$ws.Range("Q3").Value = 170
$ws.Range("R3").Value = 297
$ws.Range("S3").Value = 0.9525464349910125
$ws.Range("T3").Value = "/home/s27mhusa_hpc/Master-Thesis/Evaluation_Results/Final_TestFiles_3rdSeptember_FewShotTest_Broad/ner_evaluation_results_Qwen2.5-32B-Instruct_2_shot.txt"
$ws.Range("U3").Value = "/home/s27mhusa_hpc/Master-Thesis/Evaluation_Results/Final_TestFiles_3rdSeptember_FewShotTest_Broad/Stats/ner_evaluation_stats_Qwen2.5-32B-Instruct_2_shot.txt"
$ws.Range("V3").Value = "4 MLGPU"
$ws.Range("W3").Value = "0.065 kWh"

# X3 stays blank (empty placeholder cell), make sure a cell entry exists for it
$ws.Range("X3").NumberFormat = "@"
$ws.Range("X3").Value = " "
$ws.Range("X3").Value = ""
$ws.Range("X3").ClearFormats()
